# Fruta / hortaliza, semanal
# Prepend a new week (Primera/Segunda) of "Apio" price data at row 213,
# shifting every existing record down by 2 rows (the oldest week that falls
# off the bottom is preserved as the new final two rows of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 213-214; everything from the old row 213 onward
# (through the old row 328) shifts down to 215-330, and the sheet's used
# range grows from A1:R328 to A1:R330 automatically.
$ws.Range("A213:A214").EntireRow.Insert()

# Row 213: new "Primera" quality record for the newly-reported week.
$ws.Cells.Item(213,1).Value = 8
$ws.Cells.Item(213,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213,3).Value = "Coquimbo"
$ws.Cells.Item(213,4).Value = 44603
$ws.Cells.Item(213,5).Value = 4
$ws.Cells.Item(213,6).Value = 100112017
$ws.Cells.Item(213,7).Value = "Apio"
$ws.Cells.Item(213,8).Value = "Americana (o)"
$ws.Cells.Item(213,9).Value = "Primera"
$ws.Cells.Item(213,10).Value = 2000
$ws.Cells.Item(213,11).Value = 8500
$ws.Cells.Item(213,12).Value = 9000
$ws.Cells.Item(213,13).Value = 8750
$ws.Cells.Item(213,14).Value = "`$/docena de matas"
$ws.Cells.Item(213,15).Value = "Provincia del Elquí"
$ws.Cells.Item(213,16).Value = 1458
$ws.Cells.Item(213,17).Value = 6
$ws.Cells.Item(213,18).Value = "Hortaliza"

# Row 214: new "Segunda" quality record for the same newly-reported week.
$ws.Cells.Item(214,1).Value = 8
$ws.Cells.Item(214,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(214,3).Value = "Coquimbo"
$ws.Cells.Item(214,4).Value = 44603
$ws.Cells.Item(214,5).Value = 4
$ws.Cells.Item(214,6).Value = 100112017
$ws.Cells.Item(214,7).Value = "Apio"
$ws.Cells.Item(214,8).Value = "Americana (o)"
$ws.Cells.Item(214,9).Value = "Segunda"
$ws.Cells.Item(214,10).Value = 1320
$ws.Cells.Item(214,11).Value = 6500
$ws.Cells.Item(214,12).Value = 7000
$ws.Cells.Item(214,13).Value = 6750
$ws.Cells.Item(214,14).Value = "`$/docena de matas"
$ws.Cells.Item(214,15).Value = "Provincia del Elquí"
$ws.Cells.Item(214,16).Value = 1125
$ws.Cells.Item(214,17).Value = 6
$ws.Cells.Item(214,18).Value = "Hortaliza"
